# Team Records.xlsx - Week 13 update
# Updates each team's Record / Div. Record / Conf. Record after Week 13 games.
# Cells are touched in the same column-major order (col B..AG, then row 2,3,4
# within each column) that the original author used, so that newly introduced
# shared strings land in the same order as the authoritative workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "'7-5-0"
$ws.Range("B3").Value = "'3-1-0"
$ws.Range("B4").Value = "'5-5-0"
$ws.Range("C2").Value = "'6-7-0"
$ws.Range("D2").Value = "'3-9-0"
$ws.Range("E2").Value = "'9-4-0"
$ws.Range("E3").Value = "'3-1-0"
$ws.Range("E4").Value = "'7-1-0"
$ws.Range("F2").Value = "'7-5-0"
$ws.Range("F4").Value = "'5-3-0"
$ws.Range("H2").Value = "'8-4-0"
$ws.Range("H3").Value = "'1-2-0"
$ws.Range("H4").Value = "'5-4-0"
$ws.Range("I2").Value = "'6-5-1"
$ws.Range("I3").Value = "'2-2-0"
$ws.Range("I4").Value = "'4-4-0"
$ws.Range("J2").Value = "'7-6-0"
$ws.Range("J3").Value = "'3-2-0"
$ws.Range("J4").Value = "'6-3-0"
$ws.Range("K2").Value = "'2-10-0"
$ws.Range("L2").Value = "'2-10-0"
$ws.Range("L3").Value = "'2-2-0"
$ws.Range("L4").Value = "'2-7-0"
$ws.Range("N2").Value = "'6-6-0"
$ws.Range("N3").Value = "'1-2-0"
$ws.Range("N4").Value = "'3-5-0"
$ws.Range("O2").Value = "'7-5-0"
$ws.Range("O4").Value = "'5-3-0"
$ws.Range("P2").Value = "'8-4-0"
$ws.Range("P3").Value = "'2-1-0"
$ws.Range("P4").Value = "'3-4-0"
$ws.Range("Q2").Value = "'6-6-0"
$ws.Range("R2").Value = "'8-4-0"
$ws.Range("R4").Value = "'6-1-0"
$ws.Range("S2").Value = "'6-7-0"
$ws.Range("T2").Value = "'6-6-0"
$ws.Range("U2").Value = "'4-8-0"
$ws.Range("V2").Value = "'4-8-0"
$ws.Range("V4").Value = "'2-5-0"
$ws.Range("W2").Value = "'1-10-1"
$ws.Range("W3").Value = "'1-4-0"
$ws.Range("W4").Value = "'1-7-0"
$ws.Range("Y2").Value = "'5-7-0"
$ws.Range("Y3").Value = "'2-1-0"
$ws.Range("Y4").Value = "'4-4-0"
$ws.Range("Z2").Value = "'9-3-0"
$ws.Range("Z3").Value = "'2-1-0"
$ws.Range("Z4").Value = "'6-3-0"
$ws.Range("AA2").Value = "'5-7-0"
$ws.Range("AA3").Value = "'1-3-0"
$ws.Range("AA4").Value = "'2-6-0"
$ws.Range("AC2").Value = "'5-7-0"
$ws.Range("AC4").Value = "'5-4-0"
$ws.Range("AD2").Value = "'6-6-0"
$ws.Range("AD3").Value = "'1-4-0"
$ws.Range("AD4").Value = "'5-5-0"
$ws.Range("AE2").Value = "'10-2-0"
$ws.Range("AE4").Value = "'6-2-0"
$ws.Range("AF2").Value = "'8-4-0"
$ws.Range("AG2").Value = "'4-8-0"
$ws.Range("AG3").Value = "'2-2-0"
$ws.Range("AG4").Value = "'2-6-0"

# New column width for Lions (column W), matching the other auto-fit record
# columns once the wider "1-10-1" value is in place.
$ws.Columns("W").ColumnWidth = 9.75

# Scroll the view over and move the active selection, matching where the
# user ended up after finishing the Week 13 update (column AG, row 5).
$ws.Range("AG5").Select()
